# Updated results with new RNG
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels: re-cased to match the "new RNG" export ---
$ws.Range("C1").Value = "Integral"
$ws.Range("D1").Value = "Time"
$ws.Range("A13").Value = "Avg"

# --- New sample data from the re-run with the new RNG (B2:D11) ---
# (written as plain decimals - the PS parser here doesn't accept "E-2" exponents)
$ws.Range("B2").Value  = 0.0193272
$ws.Range("C2").Value  = 0.1962515702
$ws.Range("D2").Value  = 0.113072359

$ws.Range("B3").Value  = 0.01783849425
$ws.Range("C3").Value  = 0.1821007526
$ws.Range("D3").Value  = 0.09829286299999999

$ws.Range("B4").Value  = 0.01276449707
$ws.Range("C4").Value  = 0.14638135290000001
$ws.Range("D4").Value  = 0.071126959

$ws.Range("B5").Value  = 0.05330981606
$ws.Range("C5").Value  = 0.24086790890000001
$ws.Range("D5").Value  = 0.107232693

$ws.Range("B6").Value  = 0.20664666340000001
$ws.Range("C6").Value  = 0.37170256460000001
$ws.Range("D6").Value  = 0.103188738

$ws.Range("B7").Value  = 0.08053230607
$ws.Range("C7").Value  = 0.2723995615
$ws.Range("D7").Value  = 0.098166158

$ws.Range("B8").Value  = 0.01274157354
$ws.Range("C8").Value  = 0.15343207010000001
$ws.Range("D8").Value  = 0.102880452

$ws.Range("B9").Value  = 0.02255162186
$ws.Range("C9").Value  = 0.19250324560000001
$ws.Range("D9").Value  = 0.103446895

$ws.Range("B10").Value = 0.02111328509
$ws.Range("C10").Value = 0.19735605449999999
$ws.Range("D10").Value = 0.10368371699999999

$ws.Range("B11").Value = 0.02559282518
$ws.Range("C11").Value = 0.2119957194
$ws.Range("D11").Value = 0.10358107

# --- B13 previously only averaged B11 by mistake; fix range to match C13/D13 ---
$ws.Range("B13").Formula = "=AVERAGE(B2:B11)"
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"
$ws.Range("D13").Formula = "=AVERAGE(D2:D11)"

$ws.Range("B14").Formula = "=_xlfn.STDEV.S(B2:B11)"
$ws.Range("C14").Formula = "=_xlfn.STDEV.S(C2:C11)"
$ws.Range("D14").Formula = "=_xlfn.STDEV.S(D2:D11)"

# --- Cosmetic cleanup matching the resaved workbook ---
# Default workbook font bumped from 11pt to 12pt (resaved from a newer Excel/locale default)
$wb.Styles("Normal").Font.Size = 12

# Selection moved back to the top of the data
$ws.Range("B2").Select()

# Page margins reset to Excel's defaults
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)
